# Make the "4.6. Требования к маркировке и упаковке" heading bold,
# matching the formatting already used by the other numbered headings
# (4, 4.1 .. 4.5, 4.7) in this document. This bolds both the run text
# and the paragraph mark (so the paragraph's own rPr also gets <w:b/>).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -match "4\.6\. Требования к маркировке и упаковке") {
        $p.Range.Font.Bold = 1
        $p.Range.Font.BoldBi = 1
    }
}
